# Correction template risk indicators
# Update the uncertainty range values in columns G:P for the "block" rows
# (the first 4 rows of each 6-row group: 2-5, 8-11, 14-17, 20-23),
# replacing the narrow +/-0.1/0.05 ranges with the wider +/-0.8/0.5 ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$rowGroups = @(2, 8, 14, 20)

foreach ($base in $rowGroups) {
    $ws.Range("G" + $base + ":P" + $base).Value = -0.8
    $ws.Range("G" + ($base + 1) + ":P" + ($base + 1)).Value = -0.5
    $ws.Range("G" + ($base + 2) + ":P" + ($base + 2)).Value = 0.5
    $ws.Range("G" + ($base + 3) + ":P" + ($base + 3)).Value = 0.8
}

# Update the sheet view: scroll to column V and select cell AE6
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 22
$ws.Range("AE6").Select()
